$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B50").Value = "Lutianidae spp."
